$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the NA US address data (was a Canadian address, now a Michigan, US one)
# state (column F)
$ws.Range("F2:F5").Value = "Michigan"

# city (column E)
$ws.Range("E2:E5").Value = "TROY"

# street address (column D) - unique per row
$ws.Range("D2").Value = "973 BRAHMS CT"
$ws.Range("D3").Value = "974 BRAHMS CT"
$ws.Range("D4").Value = "975 BRAHMS CT"
$ws.Range("D5").Value = "976 BRAHMS CT"

# zip (column G) - now a plain number instead of a text postal code
$ws.Range("G2:G5").Value = 48085

# Update the sheet's saved selection
$ws.Range("A2:C5").Select()
